# Rename the embedded logo pictures that live in the document's
# first-page/default headers and footers.
#
#   - Pearson Edexcel logo (both the "first page" and "default" footers):
#         image1.png -> image2.png
#   - BTEC logo (the "first page" header):
#         image2.jpg -> image1.jpg
#
# Word COM exposes header/footer pictures as InlineShapes hanging off the
# Headers/Footers collections of a Section, so we walk those collections
# instead of $d.InlineShapes (which only covers the main document story).

$d = $word.ActiveDocument

$wdHeaderFooterPrimary   = 1
$wdHeaderFooterFirstPage = 2
$wdHeaderFooterEvenPages = 3

function Rename-LogoInlineShapes($headerFooter, [string]$newName) {
    if ($headerFooter.Exists) {
        $shapes = $headerFooter.Range.InlineShapes
        for ($i = 1; $i -le $shapes.Count; $i++) {
            $shapes.Item($i).Name = $newName
        }
    }
}

foreach ($section in $d.Sections) {

    # Pearson Edexcel logo: both footers currently named "image1.png".
    Rename-LogoInlineShapes $section.Footers.Item($wdHeaderFooterPrimary)   "image2.png"
    Rename-LogoInlineShapes $section.Footers.Item($wdHeaderFooterFirstPage) "image2.png"

    # BTEC logo lives in the first-page header, currently named "image2.jpg".
    Rename-LogoInlineShapes $section.Headers.Item($wdHeaderFooterFirstPage) "image1.jpg"
}
